$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) stays text even for values that look numeric,
# matching the workbook author's original inline-string storage.
foreach ($r in 2..51) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.149.81"
$ws.Range("D3").Value = "2.305.79"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "301.18"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "100.46"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "36.58"
$ws.Range("E10").Value = "  +8.40%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "17.78"
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("E14").Value = "  +2.53%  "
$ws.Range("D15").Value = "2.665.81"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "2.342.30"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "43.068.07"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +9.99%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  +1.90%  "
$ws.Range("D22").Value = "67.95"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "236.05"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  +8.00%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "25.16"
$ws.Range("D28").Value = "169.26"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").Value = "34.58"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "5.06"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").Value = "17.66"
$ws.Range("E34").Value = "  +5.20%  "
$ws.Range("D35").Value = "4.61"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  +1.40%  "
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("D43").Value = "1.983.64"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  -4.27%  "
$ws.Range("E45").Value = "  +4.57%  "
$ws.Range("D46").Value = "17.72"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  +2.71%  "
$ws.Range("D48").Value = "55.45"
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("E49").Value = "  +4.54%  "
$ws.Range("D50").Value = "2.531.82"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "70.85"
$ws.Range("E51").Value = "  +1.32%  "
